$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The Total Test Cases (H1) and Total Automated (H5) counts in the linked
# "_Test_Suite_Statistics_for_Folders.xlsx" source grew by one each
# (29 -> 30, 232 -> 233). Push the refreshed totals into the two cells that
# pull those external values so every dependent stat on the Update sheet
# (pass rates, SUM totals, etc.) recalculates off the new numbers.
$ws.Range("D2").Value2 = 30
$ws.Range("H2").Value2 = 233
